# Apply the "0.00" number format to the perc/baseline/UCL/LCL columns
# (B:E, rows 2-53). Excel collapses the previously-separate "General +
# applyNumberFormat" style used by column B with the "0.00" style already
# used by D/E, so B ends up sharing the same display format as C/D/E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:E53").NumberFormat = "0.00"

# Update the value (Y) axis of the scatter chart: switch it from
# "General, linked to source" to an explicit "0" integer format that is
# no longer linked to the source cell format.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)
$valueAxis.NumberFormatLinked = 0
$valueAxis.TickLabels.NumberFormat = "0"

# Move the active selection to T2.
[void]$ws.Range("T2").Select()
